$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 4 workers that left the account: delete bottom-to-top so the
# row numbers for the earlier deletes stay valid.
#   row 24 -> MARIA JOSE BROCHERO CARO  (1143378715)
#   row 22 -> MARIA JOSE MARRUGO ARAUJO (1047474987)
#   row 19 -> LINAY GALARCIO HERNANDEZ  (22810441)
#   row 17 -> ANDREINA SAAVEDRA PAYARES (1047394029)
$ws.Rows("24:24").Delete()
$ws.Rows("22:22").Delete()
$ws.Rows("19:19").Delete()
$ws.Rows("17:17").Delete()

# Update "Periodo Mora" (column E) for the remaining 8 worker rows: 2507 -> 2508
$ws.Range("E16:E23").Value = "2508"

# Update summary figures at the top of the sheet
$ws.Range("E11").Value = 1044886
$ws.Range("C13").Value = 8
